$d = $word.ActiveDocument

# --- Remove the "CryptoBridge" and "HitBTC" HYPERLINK fields ---------------
# These are plain (non-relationship) HYPERLINK fields built from fldChar /
# instrText runs. Field.Delete() removes the begin/instrText/separate/end
# runs together with the visible field-result text ("CryptoBridge",
# "HitBTC") and their trailing orphaned proofErr marks in one shot.
# (Deleting mutates $d.Fields in place, so re-scan from the top each time
# rather than continuing a stale foreach enumerator.)
$keepLooking = $true
while ($keepLooking) {
    $keepLooking = $false
    foreach ($f in $d.Fields) {
        if ($f.Code.Text -like "*crypto-bridge.org*" -or $f.Code.Text -like "*hitbtc.com*") {
            $f.Delete()
            $keepLooking = $true
            break
        }
    }
}

# --- Tidy the surrounding plain-text runs -----------------------------------
# "...exchanges such as<nbsp>" -> "...exchanges" (drop the now-dangling lead-in)
$d.Content.Find.Execute("exchanges such as" + [char]160, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "exchanges", 2)
# the ", <nbsp>" separator that used to sit between CryptoBridge and HitBTC
$d.Content.Find.Execute("," + [char]160, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# Note: saving/mutating this document also renumbers the lone "smarthosting"
# bookmark's w:id from 1 to 0, matching the target XML.
